{"js": "// The \"COMPETENCES TECHNIQUES\" skill list is being reshuffled:\n//  - the \"Bases de donn\u00e9es\" line is replaced (in place) by the\n//    \"Visualisation\" line's text, and a new \"MLOps\" line is inserted\n//    right after it;\n//  - the original \"Visualisation\" paragraph (further down) is removed;\n//  - the original \"MLOps\" paragraph (further down, now last of the\n//    group) is replaced (in place) by the \"Bases de donn\u00e9es\" text.\n// Net effect: \"Visualisation\" + \"MLOps\" move up to sit right after\n// \"Langages\", while \"Bases de donn\u00e9es\" moves down to sit right after\n// \"ML/AI\" (i.e. become the last line of the block).\n\nconst visualisationText = \"Visualisation : excel, optimization, tableau\";\nconst mlOpsText = \"MLOps : aws, spark, Git, DVC, Flask, Docker, Github Actions, Heroku, MLflow, Streamlit\";\nconst basesDeDonneesText = \"Bases de donn\u00e9es : SQL, MongoDB, Neo4j, Redis\";\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the three relevant paragraphs by their current (pre-edit) text.\nlet basesParagraph = null;\nlet visualisationParagraph = null;\nlet mlOpsParagraph = null;\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text.trim();\n  if (t === basesDeDonneesText) {\n    basesParagraph = paragraphs.items[i];\n  } else if (t === visualisationText) {\n    visualisationParagraph = paragraphs.items[i];\n  } else if (t === mlOpsText) {\n    mlOpsParagraph = paragraphs.items[i];\n  }\n}\n\nif (!basesParagraph || !visualisationParagraph || !mlOpsParagraph) {\n  throw new Error(\"Could not locate the expected skill paragraphs.\");\n}\n\n// 1) Turn the \"Bases de donn\u00e9es\" paragraph into the \"Visualisation\" one,\n//    and insert a fresh \"MLOps\" paragraph right after it.\nbasesParagraph.insertText(visualisationText, \"Replace\");\nbasesParagraph.insertParagraph(mlOpsText, \"After\");\n\n// 2) Remove the old \"Visualisation\" paragraph further down.\nvisualisationParagraph.delete();\n\n// 3) Turn the old \"MLOps\" paragraph (further down) into \"Bases de donn\u00e9es\".\nmlOpsParagraph.insertText(basesDeDonneesText, \"Replace\");\n\nawait context.sync();\n", "ps1": "# The \"COMPETENCES TECHNIQUES\" skill list is being reshuffled:\n#  - the \"Bases de donn\u00e9es\" line is replaced (in place) by the\n#    \"Visualisation\" line's text, and a new \"MLOps\" line is inserted\n#    right after it;\n#  - the original \"Visualisation\" paragraph (further down) is removed;\n#  - the original \"MLOps\" paragraph (further down, now last of the\n#    group) is replaced (in place) by the \"Bases de donn\u00e9es\" text.\n# Net effect: \"Visualisation\" + \"MLOps\" move up to sit right after\n# \"Langages\", while \"Bases de donn\u00e9es\" moves down to sit right after\n# \"ML/AI\" (i.e. become the last line of the block).\n\n$d = $word.ActiveDocument\n\n$visualisationText = \"Visualisation : excel, optimization, tableau\"\n$mlOpsText = \"MLOps : aws, spark, Git, DVC, Flask, Docker, Github Actions, Heroku, MLflow, Streamlit\"\n$basesDeDonneesText = \"Bases de donn\u00e9es : SQL, MongoDB, Neo4j, Redis\"\n\n# Locate the three relevant paragraphs by their current (pre-edit) text.\n$basesIdx = -1\n$visuIdx = -1\n$mlOpsIdx = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text.TrimEnd(\"`r\", \"`n\", \"`v\")\n    if ($t -eq $basesDeDonneesText) {\n        $basesIdx = $i\n    } elseif ($t -eq $visualisationText) {\n        $visuIdx = $i\n    } elseif ($t -eq $mlOpsText) {\n        $mlOpsIdx = $i\n    }\n}\n\nif ($basesIdx -eq -1 -or $visuIdx -eq -1 -or $mlOpsIdx -eq -1) {\n    throw \"Could not locate the expected skill paragraphs.\"\n}\n\n# 1) Turn the \"Bases de donn\u00e9es\" paragraph into the \"Visualisation\" one,\n#    and insert a fresh \"MLOps\" paragraph right after it.\n$basesParagraph = $d.Paragraphs.Item($basesIdx)\n$basesParagraph.Range.Text = $visualisationText\n$basesParagraph.Range.InsertParagraphAfter()\n$d.Paragraphs.Item($basesIdx + 1).Range.Text = $mlOpsText\n\n# 2) Remove the old \"Visualisation\" paragraph further down (indices shifted\n#    by +1 because of the paragraph inserted in step 1).\n$d.Paragraphs.Item($visuIdx + 1).Range.Delete()\n\n# 3) Turn the old \"MLOps\" paragraph (further down) into \"Bases de donn\u00e9es\".\n#    Its index shifted by +1 (insert) -1 (delete) = net 0 from the original.\n$d.Paragraphs.Item($mlOpsIdx).Range.Text = $basesDeDonneesText\n"}
